$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$req1 = "LOB1021 -  Física IV  (Requisito)`n"
$req2 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$req3 = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"

$ws.Range("B23").Value = $req1
$ws.Range("C23").Value = $req1

$ws.Range("B24").Value = $req2
$ws.Range("C24").Value = $req2

$ws.Range("B25").Value = $req3
$ws.Range("C25").Value = $req3
